# edit.ps1
# Applies the changes described by the diff:
#  - rows 175-178: fill in / update F,G,H,I,K,L,M,N,O with new order/values
#  - adds a brand new row 204 at the end of the table
#  - dimension grows from A1:O203 to A1:O204 (handled automatically by Excel
#    once data is written to row 204)
#  - two brand new shared strings ("K.02.0700" / "MATERIAL HIDRAULICO - ED")
#    get introduced through H204/I204

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: capture the "old" text values that are about to be overwritten but
# are still needed elsewhere (handles the 3-way rotation between rows
# 175/176/177 for columns H & I, plus the value that must move from N175/O175
# down to N177/O177).  Value2/Text getters round-trip plain (non numeric
# looking) text safely.
# ---------------------------------------------------------------------------
$h175old = $ws.Range("H175").Value2
$i175old = $ws.Range("I175").Value2
$h176old = $ws.Range("H176").Value2
$i176old = $ws.Range("I176").Value2
$h177old = $ws.Range("H177").Value2
$i177old = $ws.Range("I177").Value2

# ---------------------------------------------------------------------------
# Step 2: the N/O columns on rows 175-178 hold codes that look like pure
# numbers (e.g. "00000000010258"). Assigning such text through .Value would
# get auto-coerced to a number, so those moves are done with
# Copy + PasteSpecial(xlPasteValues) which preserves the text type, reuses
# the existing shared string and keeps the destination's existing style.
# Do this BEFORE overwriting N175/O175 so the old contents aren't lost.
# ---------------------------------------------------------------------------
$xlPasteValues = -4163

$ws.Range("N175").Copy() | Out-Null
$ws.Range("N177").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("O175").Copy() | Out-Null
$ws.Range("O177").PasteSpecial($xlPasteValues) | Out-Null

# N176/O176 and N178/O178 are brand new values copied from elsewhere on the
# sheet (rows untouched by this edit).
$ws.Range("N94").Copy() | Out-Null
$ws.Range("N176").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("O94").Copy() | Out-Null
$ws.Range("O176").PasteSpecial($xlPasteValues) | Out-Null

$ws.Range("N103").Copy() | Out-Null
$ws.Range("N178").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("O103").Copy() | Out-Null
$ws.Range("O178").PasteSpecial($xlPasteValues) | Out-Null

# finally overwrite N175/O175 themselves with their new value
$ws.Range("N94").Copy() | Out-Null
$ws.Range("N175").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("O94").Copy() | Out-Null
$ws.Range("O175").PasteSpecial($xlPasteValues) | Out-Null

# F176/G176, F177/G177 and F178/G178 were empty before (generic style),
# so their number format needs to be (re)applied to match F175/G175 (an
# integer doc-number column and a date column respectively) before the
# values are written. Pull the formatting from F6/G6, which already carry
# the correct styles and are untouched by this edit.
$ws.Range("F6").Copy() | Out-Null
$ws.Range("F176").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F6").Copy() | Out-Null
$ws.Range("F177").PasteSpecial(-4122) | Out-Null
$ws.Range("F6").Copy() | Out-Null
$ws.Range("F178").PasteSpecial(-4122) | Out-Null

$ws.Range("G6").Copy() | Out-Null
$ws.Range("G176").PasteSpecial(-4122) | Out-Null
$ws.Range("G6").Copy() | Out-Null
$ws.Range("G177").PasteSpecial(-4122) | Out-Null
$ws.Range("G6").Copy() | Out-Null
$ws.Range("G178").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# Step 3: row 175
# ---------------------------------------------------------------------------
$ws.Range("F175").Value = 80829
$ws.Range("G175").Value = 45985
$ws.Range("H175").Value = $h176old
$ws.Range("I175").Value = $i176old
$ws.Range("K175").Value = 3
$ws.Range("L175").Value = 80
$ws.Range("M175").Value = 240

# ---------------------------------------------------------------------------
# Step 4: row 176
# ---------------------------------------------------------------------------
$ws.Range("F176").Value = 80829
$ws.Range("G176").Value = 45985
$ws.Range("H176").Value = $h177old
$ws.Range("I176").Value = $i177old
$ws.Range("K176").Value = 1
$ws.Range("L176").Value = 275
$ws.Range("M176").Value = 275

# ---------------------------------------------------------------------------
# Step 5: row 177
# ---------------------------------------------------------------------------
$ws.Range("F177").Value = 80800
$ws.Range("G177").Value = 45982
$ws.Range("H177").Value = $h175old
$ws.Range("I177").Value = $i175old
$ws.Range("K177").Value = 12
$ws.Range("L177").Value = 13.9
$ws.Range("M177").Value = 166.8

# ---------------------------------------------------------------------------
# Step 6: row 178 (H178/I178 stay the same, only F/G/L/M/N/O change)
# ---------------------------------------------------------------------------
$ws.Range("F178").Value = 80823
$ws.Range("G178").Value = 45985
$ws.Range("L178").Value = 130
$ws.Range("M178").Value = 390

# ---------------------------------------------------------------------------
# Step 7: brand new row 204.
# Copy formatting from row 203 (the previous last row) so the new cells get
# the correct number formats/styles, then fill in the values.
# ---------------------------------------------------------------------------
$ws.Range("A203:O203").Copy() | Out-Null
$ws.Range("A204:O204").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A204").Value = 2506
$ws.Range("B124").Copy() | Out-Null
$ws.Range("B204").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C204").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("D204").Value = 23
$ws.Range("E204").Value = 45985.5804617708
$ws.Range("H204").Value = "K.02.0700"
$ws.Range("I204").Value = "MATERIAL HIDRAULICO - ED"
$ws.Range("J175").Copy() | Out-Null
$ws.Range("J204").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("K204").Value = 1
$ws.Range("L204").Value = 0
$ws.Range("M204").Value = 0

$excel.CutCopyMode = 0
